$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("I2").Value = 0.6243561829432177
$ws.Range("J2").Value = 0.7137249741492807
$ws.Range("M2").Value = 14.817691
$ws.Range("N2").Value = 29.635382
$ws.Range("O2").Value = 0.2230198506330352
$ws.Range("P2").Value = 0.1687294465231367
$ws.Range("Q2").Value = 5.465011402316666
$ws.Range("R2").Value = 32.79006841389999
$ws.Range("S2").Value = 0.1392438226618084
$ws.Range("T2").Value = 0.1204264198579482

# Row 3
$ws.Range("I3").Value = 0.6243561829432177
$ws.Range("J3").Value = 0.7137249741492807
$ws.Range("O3").Value = 0.3961505993138794
$ws.Range("P3").Value = 0.4495716715730233
$ws.Range("S3").Value = 0.2473390760582818
$ws.Range("T3").Value = 0.320870529671705

# Row 4
$ws.Range("I4").Value = 0.6243561829432177
$ws.Range("J4").Value = 0.7137249741492807
$ws.Range("M4").Value = 8.453176999999998
$ws.Range("N4").Value = 25.359531
$ws.Range("O4").Value = 0.1272280729780779
$ws.Range("P4").Value = 0.144384831270821
$ws.Range("Q4").Value = 3.117672563883332
$ws.Range("R4").Value = 28.05905307494999
$ws.Range("S4").Value = 0.07943563400781387
$ws.Range("T4").Value = 0.103051059966315

# Row 5
$ws.Range("I5").Value = 0.6243561829432177
$ws.Range("J5").Value = 0.7137249741492807
$ws.Range("M5").Value = 8.8672295
$ws.Range("N5").Value = 17.734459
$ws.Range("O5").Value = 0.133459943159757
$ws.Range("P5").Value = 0.1009713811503176
$ws.Range("Q5").Value = 3.270382026758333
$ws.Range("R5").Value = 19.62229216055
$ws.Range("S5").Value = 0.0833265406870447
$ws.Range("T5").Value = 0.07206579640132757

# Row 6
$ws.Range("I6").Value = 0.6243561829432177
$ws.Range("J6").Value = 0.7137249741492807
$ws.Range("M6").Value = 3.775572333333333
$ws.Range("N6").Value = 11.326717
$ws.Range("O6").Value = 0.05682582919526532
$ws.Range("P6").Value = 0.06448881577886201
$ws.Range("Q6").Value = 1.392494002738889
$ws.Range("R6").Value = 12.53244602465
$ws.Range("S6").Value = 0.03547955780893912
$ws.Range("T6").Value = 0.04602727837468601

# Row 7
$ws.Range("I7").Value = 0.6243561829432177
$ws.Range("J7").Value = 0.7137249741492807
$ws.Range("M7").Value = 4.206767
$ws.Range("N7").Value = 12.620301
$ws.Range("O7").Value = 0.06331570471998517
$ws.Range("P7").Value = 0.07185385370383916
$ws.Range("Q7").Value = 1.551525782383333
$ws.Range("R7").Value = 13.96373204145
$ws.Range("S7").Value = 0.03953155171932982
$ws.Range("T7").Value = 0.0512838898772988

# Row 8
$ws.Range("E8").Value = 1
$ws.Range("F8").Value = 0.5
$ws.Range("G8").Value = 0.2218985
$ws.Range("H8").Value = 0.443797
$ws.Range("I8").Value = 0.3756438170567823
$ws.Range("J8").Value = 0.2862750258507193
$ws.Range("M8").Value = 14.817691
$ws.Range("N8").Value = 29.635382
$ws.Range("O8").Value = 0.2230198506330352
$ws.Range("P8").Value = 0.1687294465231367
$ws.Range("Q8").Value = 3.2880234063635
$ws.Range("R8").Value = 13.152093625454
$ws.Range("S8").Value = 0.0837760279712268
$ws.Range("T8").Value = 0.04830302666518851

# Row 9
$ws.Range("E9").Value = 1
$ws.Range("F9").Value = 0.5
$ws.Range("G9").Value = 0.2218985
$ws.Range("H9").Value = 0.443797
$ws.Range("I9").Value = 0.3756438170567823
$ws.Range("J9").Value = 0.2862750258507193
$ws.Range("O9").Value = 0.3961505993138794
$ws.Range("P9").Value = 0.4495716715730233
$ws.Range("Q9").Value = 5.840522443592834
$ws.Range("R9").Value = 35.043134661557
$ws.Range("S9").Value = 0.1488115232555976
$ws.Range("T9").Value = 0.1287011419013183

# Row 10
$ws.Range("E10").Value = 1
$ws.Range("F10").Value = 0.5
$ws.Range("G10").Value = 0.2218985
$ws.Range("H10").Value = 0.443797
$ws.Range("I10").Value = 0.3756438170567823
$ws.Range("J10").Value = 0.2862750258507193
$ws.Range("M10").Value = 8.453176999999998
$ws.Range("N10").Value = 25.359531
$ws.Range("O10").Value = 0.1272280729780779
$ws.Range("P10").Value = 0.144384831270821
$ws.Range("Q10").Value = 1.8757472965345
$ws.Range("R10").Value = 11.254483779207
$ws.Range("S10").Value = 0.04779243897026406
$ws.Range("T10").Value = 0.04133377130450603

# Row 11
$ws.Range("E11").Value = 1
$ws.Range("F11").Value = 0.5
$ws.Range("G11").Value = 0.2218985
$ws.Range("H11").Value = 0.443797
$ws.Range("I11").Value = 0.3756438170567823
$ws.Range("J11").Value = 0.2862750258507193
$ws.Range("M11").Value = 8.8672295
$ws.Range("N11").Value = 17.734459
$ws.Range("O11").Value = 0.133459943159757
$ws.Range("P11").Value = 0.1009713811503176
$ws.Range("Q11").Value = 1.96762492520575
$ws.Range("R11").Value = 7.870499700823
$ws.Range("S11").Value = 0.05013340247271234
$ws.Range("T11").Value = 0.02890558474898999

# Row 12
$ws.Range("E12").Value = 1
$ws.Range("F12").Value = 0.5
$ws.Range("G12").Value = 0.2218985
$ws.Range("H12").Value = 0.443797
$ws.Range("I12").Value = 0.3756438170567823
$ws.Range("J12").Value = 0.2862750258507193
$ws.Range("M12").Value = 3.775572333333333
$ws.Range("N12").Value = 11.326717
$ws.Range("O12").Value = 0.05682582919526532
$ws.Range("P12").Value = 0.06448881577886201
$ws.Range("Q12").Value = 0.8377938374081666
$ws.Range("R12").Value = 5.026763024448999
$ws.Range("S12").Value = 0.02134627138632621
$ws.Range("T12").Value = 0.01846153740417599

# Row 13
$ws.Range("E13").Value = 1
$ws.Range("F13").Value = 0.5
$ws.Range("G13").Value = 0.2218985
$ws.Range("H13").Value = 0.443797
$ws.Range("I13").Value = 0.3756438170567823
$ws.Range("J13").Value = 0.2862750258507193
$ws.Range("M13").Value = 4.206767
$ws.Range("N13").Value = 12.620301
$ws.Range("O13").Value = 0.06331570471998517
$ws.Range("P13").Value = 0.07185385370383916
$ws.Range("Q13").Value = 0.9334752871495
$ws.Range("R13").Value = 5.600851722897001
$ws.Range("S13").Value = 0.02378415300065536
$ws.Range("T13").Value = 0.02056996382654036
